# Actualización automática 2025-09-29 16:30:10
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "VENTAS POR GRUPO" (sheet1)
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Range("I51").Value = 99
$ws1.Range("N51").Value = 124.22
$ws1.Range("O51").Value = 963.53

$ws1.Range("N54").Value = "2 de 52"
$ws1.Range("O54").Value = "4 de 52"

# ---------------------------------------------------------------
# Sheet "VENTA MENSUAL" (sheet2)
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Range("F53").Value = 1709.75
$ws2.Range("F54").Value = 1709.75
$ws2.Range("F58").Value = 75129.67

# ---------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL" (sheet3)
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Column F width 25 -> 24
# (Excel's ColumnWidth<->stored-width conversion rounds via pixels, so use a
#  value that round-trips to a stored width of exactly 24.)
$ws3.Columns.Item(6).ColumnWidth = 23.16

$ws3.Range("D7").Value = 2396.7
$ws3.Range("E7").Value = -1509.988983712426
$ws3.Range("F7").Value = 2.702909917635119

$ws3.Range("D8").Value = 0
$ws3.Range("E8").Value = 1346.40488751609
$ws3.Range("F8").Value = 0

$ws3.Range("D13").Value = 358.6
$ws3.Range("E13").Value = 82.05317777811899
$ws3.Range("F13").Value = 0.81379192998936

$ws3.Range("D14").Value = 4595.09
$ws3.Range("E14").Value = 3242.22410570622
$ws3.Range("F14").Value = 0.5863092812184713

$ws3.Range("D15").Value = 72921.91
$ws3.Range("E15").Value = 49132.92551083436
$ws3.Range("F15").Value = 0.5974520361671947
